# Added save file flag
# Update the row-2 results with the refreshed simulation output values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 337.60000000000002
$ws.Range("B2").Value = 290
$ws.Range("C2").Value = 16.095753678750174
$ws.Range("D2").Value = 34.879336484153207
$ws.Range("E2").Value = 15.703965597336769
$ws.Range("F2").Value = 30.860472118130232

# The shorter values above render narrower, so the sheet's columns were
# resized (re-autofit) to hug the new content, matching the saved file.
$ws.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws.Columns.Item(3).ColumnWidth = 18
$ws.Columns.Item(4).ColumnWidth = 20.333333333333332
$ws.Columns.Item(5).ColumnWidth = 16.666666666666668
$ws.Columns.Item(6).ColumnWidth = 18.833333333333332
